$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: "Pathfinding Package added. Objects can now be thrown."
# Date: 2017-09-08 (serial 42986, same day as row 7)
# D8: description text (wrap-text style like D7)
# F8/G8: two reference links (plain style like F5/G5)

# Copy formatting from analogous existing cells so the same shared
# cellXfs (styles) are reused instead of new ones being created.
# New shared-string entries are created in the same order they appear
# in the target file (astar link, ironic.games link, then the Nav-Mesh
# description).
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A8").Value = 42986

$ws.Range("F5").Copy()
$ws.Range("F8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F8").Value = "https://arongranberg.com/astar/download"

$ws.Range("G5").Copy()
$ws.Range("G8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G8").Value = "https://ironic.games/unity-asset-store-reviews/unity-navmesh-vs-a-star-pathfinding-vs-apex-path"

$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D8").Value = "Nav-Mesh: mit Unity nicht in 2D, ohne weiteres"

$ws.Rows.Item(8).RowHeight = 30

# Update the active selection, matching the saved view state.
$ws.Range("D9").Select() | Out-Null
